# All_tc_results.xlsx edit
#
# Commit message: "removed superflouous text from scripts, added dimension
# pass fail data for tc4 and 5"
#
# Net effect (per the OOXML diff):
#   - Test_Case_4 and Test_Case_5 sheets: the "dimension" comparison block
#     (rows 2-5, columns D:G) used to hold a handful of raw numbers and a
#     lot of truly-empty string cells. It now holds a full 4-row x 4-column
#     table of dimension pass/fail numbers (stored as TEXT, matching the
#     rest of that block which was already text), plus row 5 is now
#     populated with "Pass".
#   - Row 4 (the "Percent Difference" row) values B4/C4 changed from
#     0.0769...  to 7.69...  (values were rescaled from fraction to percent).
#   - Test_Case_5 row 3 B3/C3 (Modpath6 totals) were corrected to new
#     numbers.
#
# Because all of these written values look numeric (e.g. "2222.0106..."),
# a plain `.Value = "2222.0106..."` assignment would make Excel infer a
# Number cell - but the target file stores them as shared-string TEXT
# cells (t="s"), matching how the rest of that block was authored. We
# force text entry the same way Excel's UI does for a typed numeric
# string: a leading apostrophe. That marks the cell "quote-prefixed"
# (a text-display style) so we immediately reset the cell style back to
# "Normal" afterwards - the stored value stays text, but no visible
# formatting sticks around on the cell.

function Set-TextValue($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Test_Case_4
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Test_Case_4")

$tc4 = [ordered]@{
    "D2" = "2222.010614165105"
    "E2" = "1952.2631030716002"
    "F2" = "6013.278330701403"
    "G2" = "6507.326095201075"

    "D3" = "2333.68519866677"
    "E3" = "2220.58759854444"
    "F3" = "5900.079250175981"
    "G3" = "6088.435882829661"

    "B4" = "7.690028701551169"
    "C4" = "7.690028701551169"
    "D4" = "2.4513178467077372"
    "E4" = "6.430244325993366"
    "F4" = "0.9501862070112219"
    "G4" = "3.325644078556214"

    "D5" = "Pass"
    "E5" = "Pass"
    "F5" = "Pass"
    "G5" = "Pass"
}

foreach ($addr in $tc4.Keys) {
    Set-TextValue $ws4 $addr $tc4[$addr]
}

# ---------------------------------------------------------------------
# Test_Case_5
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Test_Case_5")

$tc5 = [ordered]@{
    "D2" = "2222.010614165105"
    "E2" = "1952.2631030716002"
    "F2" = "6013.278330701403"
    "G2" = "6507.326095201075"

    "B3" = "9349386.234604424"
    "C3" = "214.63385978781378"
    "D3" = "2332.568630426191"
    "E3" = "2061.142585285008"
    "F3" = "5900.046642815694"
    "G3" = "5929.750451251864"

    "B4" = "7.743466775023054"
    "C4" = "7.743466775023054"
    "D4" = "2.4274035058754757"
    "E4" = "2.712894999109628"
    "F4" = "0.9504625126689582"
    "G4" = "4.6439823843806405"

    "D5" = "Pass"
    "E5" = "Pass"
    "F5" = "Pass"
    "G5" = "Pass"
}

foreach ($addr in $tc5.Keys) {
    Set-TextValue $ws5 $addr $tc5[$addr]
}
